# Updates cryptos list (Price / Volume(1h) columns, plus one coin swap in row 51)
# to match the latest scrape, per the "Updated cryptos list ... with GitHub Actions" commit.
#
# Note: Price values in column D are stored as text in the workbook (e.g. "1.002",
# "20.67"), not numbers, even when they look numeric. Assigning a bare numeric-looking
# string via .Value would make Excel auto-convert it to a real number, so those are
# written with a leading apostrophe ('1.002) to force text entry, exactly like typing
# it into Excel's UI. Values that aren't numeric-looking (e.g. "27.482.58",
# "1.803.57") are unaffected by this and are assigned directly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.482.58'
$ws.Range('E2').Value = '  -0.81%  '
$ws.Range('D3').Value = '1.830.76'
$ws.Range('E3').Value = '  -0.88%  '
$ws.Range('D4').Value = "'1.002"
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = "'312.52"
$ws.Range('E5').Value = '  -0.45%  '
$ws.Range('E6').Value = '  +0.04%  '
$ws.Range('D7').Value = "'0.4289"
$ws.Range('E7').Value = '  -0.21%  '
$ws.Range('D8').Value = "'0.3657"
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').Value = "'0.07291"
$ws.Range('E9').Value = '  -0.71%  '
$ws.Range('D10').Value = "'0.8655"
$ws.Range('E10').Value = '  -1.27%  '
$ws.Range('D11').Value = "'20.67"
$ws.Range('E11').Value = '  -0.33%  '
$ws.Range('D12').Value = '1.803.57'
$ws.Range('E12').Value = '  -1.02%  '
$ws.Range('D13').Value = "'5.406"
$ws.Range('E13').Value = '  +1.12%  '
$ws.Range('D14').Value = "'6.522"
$ws.Range('E14').Value = '  -0.08%  '
$ws.Range('E15').Value = '  -0.33%  '
$ws.Range('D16').Value = "'1.003"
$ws.Range('E16').Value = '  +0.01%  '
$ws.Range('D17').Value = "'80.48"
$ws.Range('E17').Value = '  +0.83%  '
$ws.Range('D18').Value = "'0.000008922"
$ws.Range('E18').Value = '  -1.06%  '
$ws.Range('D19').Value = "'1.002"
$ws.Range('E19').Value = '  +0.13%  '
$ws.Range('D20').Value = "'15.41"
$ws.Range('E20').Value = '  +0.15%  '
$ws.Range('D21').Value = '27.573.77'
$ws.Range('E21').Value = '  -0.01%  '
$ws.Range('D22').Value = "'5.135"
$ws.Range('E22').Value = '  +3.35%  '
$ws.Range('D23').Value = "'10.84"
$ws.Range('E23').Value = '  +4.68%  '
$ws.Range('D24').Value = '2.054.33'
$ws.Range('E24').Value = '  -1.43%  '
$ws.Range('D25').Value = "'1.982"
$ws.Range('E25').Value = '  +0.18%  '
$ws.Range('D26').Value = "'154.38"
$ws.Range('E26').Value = '  -1.17%  '
$ws.Range('D27').Value = "'18.99"
$ws.Range('E27').Value = '  +1.84%  '
$ws.Range('D28').Value = "'5.126"
$ws.Range('E28').Value = '  -2.38%  '
$ws.Range('D29').Value = "'114.45"
$ws.Range('E29').Value = '  -3.98%  '
$ws.Range('D30').Value = "'1.836"
$ws.Range('E30').Value = '  -2.06%  '
$ws.Range('D31').Value = "'0.08869"
$ws.Range('E31').Value = '  -0.42%  '
$ws.Range('D32').Value = "'0.7535"
$ws.Range('E32').Value = '  +0.03%  '
$ws.Range('D33').Value = "'2.980"
$ws.Range('E33').Value = '  +0.42%  '
$ws.Range('D34').Value = "'4.539"
$ws.Range('E34').Value = '  -0.26%  '
$ws.Range('D35').Value = "'1.136"
$ws.Range('E35').Value = '  +1.29%  '
$ws.Range('E36').Value = '  +0.09%  '
$ws.Range('E37').Value = '  -1.50%  '
$ws.Range('D38').Value = "'0.05325"
$ws.Range('E38').Value = '  -2.08%  '
$ws.Range('D39').Value = "'0.01939"
$ws.Range('E39').Value = '  +0.28%  '
$ws.Range('D40').Value = "'2.800"
$ws.Range('E40').Value = '  -0.88%  '
$ws.Range('D41').Value = "'0.1666"
$ws.Range('E41').Value = '  +0.39%  '
$ws.Range('D42').Value = "'0.5088"
$ws.Range('E42').Value = '  +0.01%  '
$ws.Range('D43').Value = "'6.582"
$ws.Range('E43').Value = '  -0.03%  '
$ws.Range('D44').Value = "'8.343"
$ws.Range('E44').Value = '  -0.24%  '
$ws.Range('D45').Value = "'10.46"
$ws.Range('E45').Value = '  +1.07%  '
$ws.Range('D46').Value = "'106.03"
$ws.Range('E46').Value = '  +0.82%  '
$ws.Range('D47').Value = "'0.06499"
$ws.Range('E47').Value = '  -0.69%  '
$ws.Range('D48').Value = "'0.4686"
$ws.Range('E48').Value = '  +0.39%  '
$ws.Range('E49').Value = '  +0.04%  '
$ws.Range('D50').Value = "'1.617"
$ws.Range('E50').Value = '  -0.79%  '
$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D51').Value = "'63.94"
$ws.Range('E51').Value = '  -0.62%  '

Write-Output "Applied cryptos update"
